$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.756.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.33%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.937.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.32%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'352.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.96%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'105.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -4.52%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E8").Value = "'  -0.02%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -5.64%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'37.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -5.09%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +2.03%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.0845"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.86%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'18.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -4.60%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'3.402.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.34%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'7.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -5.91%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.935.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.08%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.976"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.31%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'51.667.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.59%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -1.28%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -4.15%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'13.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -6.37%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -3.03%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "'BitcoinCash"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'265.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.77%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "'Litecoin"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'68.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.23%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -6.46%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -6.74%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'26.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.37%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.10%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -5.51%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.107"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.42%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +2.73%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -5.61%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -5.28%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'35.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -7.63%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'50.83"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.51%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.0426"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -4.41%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.05%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -1.34%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +3.41%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -6.48%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -5.65%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -4.27%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'23.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.62%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'120.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.08%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -0.73%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.097.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.72%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -7.60%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -7.29%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'3.230.43"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.36%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -5.15%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -5.10%  "
$ws.Range("E51").Style = "Normal"
